$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 211, shifting existing rows 211+ down by one.
$ws.Rows.Item(211).Insert()

# Populate the newly inserted row 211 with the new data record.
$ws.Cells.Item(211, 1).Value = 7
$ws.Cells.Item(211, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(211, 3).Value = "Ñuble"
$ws.Cells.Item(211, 4).Value = 45141
$ws.Cells.Item(211, 5).Value = 16
$ws.Cells.Item(211, 6).Value = 100112040
$ws.Cells.Item(211, 7).Value = "Cilantro"
$ws.Cells.Item(211, 8).Value = "Sin especificar"
$ws.Cells.Item(211, 9).Value = "Primera"
$ws.Cells.Item(211, 10).Value = 60
$ws.Cells.Item(211, 11).Value = 1500
$ws.Cells.Item(211, 12).Value = 1500
$ws.Cells.Item(211, 13).Value = 1500
$ws.Cells.Item(211, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(211, 15).Value = "Región de Ñuble"
$ws.Cells.Item(211, 16).Value = 1500
$ws.Cells.Item(211, 17).Value = 1
$ws.Cells.Item(211, 18).Value = "Hortaliza"
